$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the style/format of the
# existing header cells (copy G1's formatting onto H1 via PasteSpecial)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the value for the new "Save" column in row 2
$ws.Range("H2").Value = 1
